# Weekly update for Fruta/Hortalizas - Vega Central Mapocho de Santiago - Uva
# Inserts two new daily price rows at the top of the data block (row 368),
# pushing the existing rows down by two and growing the used range to A1:T396.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 368 (shifts old rows 368:394 down to 370:396)
$ws.Rows("368:369").Insert()

# ---- New row 368 ----
$ws.Cells.Item(368, 1).Value = 9
$ws.Cells.Item(368, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(368, 3).Value = "Metropolitana"
$ws.Cells.Item(368, 4).Value = 44449
$ws.Cells.Item(368, 5).Value = 13
$ws.Cells.Item(368, 6).Value = "Fruta"
$ws.Cells.Item(368, 7).Value = 100109
$ws.Cells.Item(368, 8).Value = "Uva"
$ws.Cells.Item(368, 9).Value = 100109001
$ws.Cells.Item(368, 10).Value = "Uva"
$ws.Cells.Item(368, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(368, 12).Value = "Segunda"
$ws.Cells.Item(368, 13).Value = 35
$ws.Cells.Item(368, 14).Value = 12000
$ws.Cells.Item(368, 15).Value = 12000
$ws.Cells.Item(368, 16).Value = 12000
$ws.Cells.Item(368, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(368, 18).Value = "Región Metropolitana"
$ws.Cells.Item(368, 19).Value = 1500
$ws.Cells.Item(368, 20).Value = 8

# ---- New row 369 ----
$ws.Cells.Item(369, 1).Value = 9
$ws.Cells.Item(369, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(369, 3).Value = "Metropolitana"
$ws.Cells.Item(369, 4).Value = 44449
$ws.Cells.Item(369, 5).Value = 13
$ws.Cells.Item(369, 6).Value = "Fruta"
$ws.Cells.Item(369, 7).Value = 100109
$ws.Cells.Item(369, 8).Value = "Uva"
$ws.Cells.Item(369, 9).Value = 100109001
$ws.Cells.Item(369, 10).Value = "Uva"
$ws.Cells.Item(369, 11).Value = "Red Globe"
$ws.Cells.Item(369, 12).Value = "Primera"
$ws.Cells.Item(369, 13).Value = 25
$ws.Cells.Item(369, 14).Value = 8000
$ws.Cells.Item(369, 15).Value = 8000
$ws.Cells.Item(369, 16).Value = 8000
$ws.Cells.Item(369, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(369, 18).Value = "Región Metropolitana"
$ws.Cells.Item(369, 19).Value = 800
$ws.Cells.Item(369, 20).Value = 10
